$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.172.56'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.65%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.428.61'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.45%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '406.92'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.67%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '132.43'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.21%  '
$ws.Range("E7").Value = '  -2.06%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  +2.65%  '
$ws.Range("E10").Value = '  +6.89%  '
$ws.Range("E11").Value = '  -0.31%  '
$ws.Range("E12").Value = '  -0.12%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.84'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.06%  '
$ws.Range("E14").Value = '  -1.40%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.449.10'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.63%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '11.65'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.00%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.184.08'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.98%  '
$ws.Range("E18").Value = '  -0.70%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0000147'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +9.25%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.16'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.46%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '83.92'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.43%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '311.78'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.66%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.77'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.82%  '
$ws.Range("E24").Value = '  +0.80%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.71'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.10%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '29.66'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.25%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.16'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.76%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.79'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.54%  '
$ws.Range("E29").Value = '  +5.78%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.172'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.64%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '43.69'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.78%  '
$ws.Range("E32").Value = '  -0.65%  '
$ws.Range("E33").Value = '  -3.66%  '
$ws.Range("E34").Value = '  -0.05%  '
$ws.Range("E35").Value = '  +0.32%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '51.75'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.98%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.998'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.33%  '
$ws.Range("E38").Value = '  +0.97%  '
$ws.Range("E39").Value = '  -3.06%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.318'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +11.75%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '144.25'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.22%  '
$ws.Range("E42").Value = '  -0.55%  '
$ws.Range("B44").Value = 'Celestia'
$ws.Range("C44").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.81'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.82%  '
$ws.Range("B45").Value = 'NEARProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.90'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.42%  '
$ws.Range("E46").Value = '  +0.06%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '21.18'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.77%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.103.61'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.64%  '
$ws.Range("E49").Value = '  -1.88%  '
$ws.Range("E50").Value = '  +2.67%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.72'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +18.82%  '
